# Refine the June 28 "IG: Monday June 28" opening slide (slide 3):
#  - grow/shift the body placeholder
#  - expand the ITU-T bullet (new quarter refs + second presenter)
#  - expand the ISO TC184/SC4 bullet (list of working groups)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(2)        # "Content Placeholder 2"
$tr = $shp.TextFrame.TextRange

$dash = [char]0x2013            # en dash "-" used throughout the bullet list

# --- 1. Reposition / resize the placeholder -----------------------------
$shp.Top = 1570892 / 12700
$shp.Height = 4606071 / 12700

# --- 2. "ITU-T (20m) - Gyu Myoung Lee" bullet ---------------------------
# Locate the paragraph so later edits keep working even if text shifts.
$paraCount = $tr.Paragraphs().Count
$ituParaIdx = -1
for ($i = 1; $i -le $paraCount; $i++) {
    if ($tr.Paragraphs($i).Text.StartsWith("ITU-T")) {
        $ituParaIdx = $i
        break
    }
}
$ituPara = $tr.Paragraphs($ituParaIdx)

# 2a. "ITU-T (20m) - " -> "ITU-T Q2/20 & Q4/20 (20m) - "
$oldHead = "ITU-T (20m) " + $dash + " "
$headRange = $tr.Characters($ituPara.Start, $oldHead.Length)
$headRange.Text = "ITU-T Q2/20 & Q4/20 (20m) " + $dash + " "

# 2b. " Lee" -> " Lee/Marco " + new run "Carugi"
$ituPara = $tr.Paragraphs($ituParaIdx)
$relIdx = $ituPara.Text.IndexOf(" Lee")
$leeStart = $ituPara.Start + $relIdx
$leeRange = $tr.Characters($leeStart, 4)
$leeRange.Text = " Lee/Marco Carugi"

# Split "Carugi" into its own run (mirrors how the other presenter names
# on this slide - Gyu / Myoung / Sonoda / Jackle - are each their own run).
$ituPara = $tr.Paragraphs($ituParaIdx)
$relIdx2 = $ituPara.Text.IndexOf("Carugi")
$carugiStart = $ituPara.Start + $relIdx2
$carugiRange = $tr.Characters($carugiStart, 6)
$carugiRange.Text = "Carugi"

# --- 3. ISO TC184/SC4 bullet --------------------------------------------
$paraCount = $tr.Paragraphs().Count
$isoParaIdx = -1
for ($i = 1; $i -le $paraCount; $i++) {
    if ($tr.Paragraphs($i).Text.StartsWith("ISO TC184")) {
        $isoParaIdx = $i
        break
    }
}
$isoPara = $tr.Paragraphs($isoParaIdx)

$oldIso = "ISO TC184/SC4 (about IEC CDD) (20m) " + $dash + " Hiroshi Murayama/Yoshiaki "
$isoRange = $tr.Characters($isoPara.Start, $oldIso.Length)
$isoRange.Text = "ISO TC184/SC4 (incl. JWG24, IEC SC3D, WG12, WG23, OTD, IEC CDD, etc.) (20m) " + $dash + " Hiroshi Murayama/Yoshiaki "
